# RBA 2.3 - Relatório e Email
$d = $word.ActiveDocument

# 1) Main body text: "QWREW" -> "QWR" (bold run before the comma)
$d.Content.Find.Execute("QWREW", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QWR", 2)

# 2) Header text replacements
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range

# "DIRETORIA DE ENSINO REGIAO REW" -> "...QWER"
$hdrRange.Find.Execute("REW", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "QWER", 2)

# "QWREW" -> "QWR" (before " - DEP.")
$hdrRange.Find.Execute("QWREW", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "QWR", 2)

# "Rew" -> "Qwer" (5 occurrences in address line)
for ($i = 0; $i -lt 5; $i++) {
    $hdrRange.Find.Execute("Rew", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Qwer", 2)
}

# "rew" -> "qwer" (3 occurrences: CEP, Tel, Email lines)
for ($i = 0; $i -lt 3; $i++) {
    $hdrRange.Find.Execute("rew", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "qwer", 2)
}
